# Edit: swap the "Integral / Red Violet" theme colours for the
# "Office Theme / Office" colours (the deck's two theme parts had their
# content swapped), and switch the three balance-sheet tables (slides
# 14-16) from the custom "Table_0" style over to the built-in
# "No Style, No Grid" table style.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Theme colour scheme: Red Violet -> Office
#    (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
#    NB: Table.Style / ThemeColorScheme.Colors(i).RGB use the classic
#    OLE COLORREF byte order (0x00BBGGRR), not 0x00RRGGBB.
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000  # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444  # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED  # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244  # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70  # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305  # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95  # folHlink -> 954F72

# ---------------------------------------------------------------
# 2) Table styles: Table_0 {9F9C856A-...} -> No Style, No Grid
#    {23FF835A-9F1F-411B-9624-F5B935370B19}, for the tables on
#    slides 14, 15 and 16 (each slide's first shape is the table).
# ---------------------------------------------------------------
$newStyle = "{23FF835A-9F1F-411B-9624-F5B935370B19}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyle)
    }
}
